$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name (reflects the new "through" date)
$ws.Name = "Through 2021-12-30"

# Update the December row label to reflect the new date
$ws.Range("A13").Value = "December (through 12-30)"

# Update December row (row 13) values
$ws.Range("B13").Value = 46
$ws.Range("C13").Value = 96
$ws.Range("D13").Value = 114
$ws.Range("E13").Value = 75
$ws.Range("F13").Value = 64
$ws.Range("G13").Value = 145
$ws.Range("H13").Value = 192

# Update Total row (row 14) values
$ws.Range("B14").Value = 337
$ws.Range("C14").Value = 659
$ws.Range("D14").Value = 935
$ws.Range("E14").Value = 757
$ws.Range("F14").Value = 598
$ws.Range("G14").Value = 1409
$ws.Range("H14").Value = 1835
